$wb = $excel.ActiveWorkbook

# Update the "createUser" sheet: change the test user numeric suffix from 162 to 1003.
# Dependent formulas in B2 and F2 (which build usernames/emails from A2) recalc automatically.
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1003

# Update the "addListItem" sheet: change the list item name from UserEscC to UserEscD.
# C2 formula (=A2) recalculates automatically to reflect the new value.
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "UserEscD"

# Move the selected/active tab from "createUser" to "addListItem".
$wsAddListItem.Activate()
